# Auto-generated script to update cryptos.xlsx Price (D) and Volume(1h) (E) columns
# to match the target OOXML diff (Thu Apr 27 17:40:38 UTC 2023 GitHub Actions update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.265.12"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -1.87%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.904.20"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -2.63%  "
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9989"
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.12%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "333.26"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -2.62%  "
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9992"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4625"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -3.25%  "
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4139"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -0.25%  "
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "47.79"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -0.31%  "
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.08031"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -2.82%  "
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.017"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -2.14%  "
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "22.20"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -2.62%  "
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.856.29"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -4.44%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.950"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -3.71%  "
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -3.82%  "
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "89.19"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -3.17%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.9986"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.00001031"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -2.62%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06568"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -1.78%  "
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.66"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -1.96%  "
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -0.12%  "
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "29.208.62"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -1.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.509"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -1.50%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.41"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +1.09%  "
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -3.75%  "
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.154.84"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -1.09%  "
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "157.27"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -2.56%  "
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "19.80"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -2.32%  "
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.139"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -2.02%  "
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.667"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "117.22"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -4.71%  "
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.046"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +3.84%  "
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.09437"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -2.31%  "
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.427"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -3.68%  "
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.546"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -3.79%  "
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.374"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -2.68%  "
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.06107"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -2.86%  "
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02257"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -2.55%  "
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "8.450"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -0.31%  "
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.178"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -0.59%  "
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.5875"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -3.71%  "
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9987"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1830"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -3.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "10.18"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -5.40%  "
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.251"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -2.32%  "
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.347"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -2.45%  "
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.07508"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +1.80%  "
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.5566"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -2.52%  "
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "12.17"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -2.95%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.927"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -3.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "113.36"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = "Normal"
